# Changes in Actions and Scripts. IndhuTestData Added
#
# Add a new row of sign-in test data (email + password) to the "SignIn"
# sheet, resize column B to fit the new value, then switch the active
# sheet/tab over to "SignIn" with cell A6 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SignIn")

# New row 3: test credentials (IndhuTestData)
$ws.Range("A3").Value = "oct21@yopmail.com"
$ws.Range("B3").Value = "Hiindhu11@"

# Column B (Password) widens to fit the new, longer value
$ws.Columns.Item(2).ColumnWidth = 10

# Make "SignIn" the active/selected sheet and tab
$ws.Activate()

# Update the active selection on the SignIn sheet
$ws.Range("A6").Select() | Out-Null
